$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 125000000
$ws.Range("J6").Value = 125000000
$ws.Range("L6").Value = 375000000
$ws.Range("N6").Value = -375000224

$ws.Range("H31").Value = 7626.8335
$ws.Range("I31").Value = 7437.8
$ws.Range("K31").Value = 22313.4
$ws.Range("M31").Value = -22083.4

$ws.Range("H70").Value = 3226.3333
$ws.Range("J70").Value = 3338.077
$ws.Range("L70").Value = 10014.231
$ws.Range("N70").Value = -10554.231

$ws.Range("H73").Value = 3226.3333
$ws.Range("J73").Value = 3338.077
$ws.Range("L73").Value = 10014.231
$ws.Range("N73").Value = -11886.231

$ws.Range("H96").Value = 652.1667
$ws.Range("I96").Value = 652.1667
$ws.Range("K96").Value = 1956.5001
$ws.Range("M96").Value = -583.5001

$ws.Range("H100").Value = 1405.4615
$ws.Range("I100").Value = 1355.4166
$ws.Range("K100").Value = 1355.4166
$ws.Range("M100").Value = -814.4166

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 698.2857
$ws.Range("J5").Value = 698.2857
$ws.Range("L5").Value = 698.2857
$ws.Range("N5").Value = -922.2857

$ws.Range("H32").Value = 10847.1
$ws.Range("I32").Value = 9421.704
$ws.Range("K32").Value = 9421.704
$ws.Range("M32").Value = -9134.704

$ws.Range("H45").Value = 109780.1
$ws.Range("I45").Value = 184937.19
$ws.Range("J45").Value = 6439.125
$ws.Range("K45").Value = 184937.19
$ws.Range("L45").Value = 6439.125
$ws.Range("M45").Value = -184560.19
$ws.Range("N45").Value = -7193.125

$ws.Range("H64").Value = 22101
$ws.Range("J64").Value = 22101
$ws.Range("L64").Value = 22101
$ws.Range("N64").Value = -22597

$ws.Range("H67").Value = 22101
$ws.Range("J67").Value = 22101
$ws.Range("L67").Value = 22101
$ws.Range("N67").Value = -23817

$ws.Range("H102").Value = 2818.7827
$ws.Range("I102").Value = 1944.381
$ws.Range("K102").Value = 1944.381
$ws.Range("M102").Value = -322.3810000000001

$ws.Range("H110").Value = 4125.7896
$ws.Range("I110").Value = 2879.3076
$ws.Range("K110").Value = 2879.3076
$ws.Range("M110").Value = -834.3076000000001

$ws.Range("H122").Value = 1435.3334
$ws.Range("I122").Value = 1322.4
$ws.Range("K122").Value = 3967.2
$ws.Range("M122").Value = -1517.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 698.2857
$ws.Range("J4").Value = 698.2857
$ws.Range("L4").Value = 698.2857
$ws.Range("N4").Value = -928.2857

$ws.Range("H94").Value = 2334.8
$ws.Range("J94").Value = 2514.2856
$ws.Range("L94").Value = 2514.2856
$ws.Range("N94").Value = -3416.2856

$ws.Range("H99").Value = 4910.1816
$ws.Range("I99").Value = 3963.4614
$ws.Range("J99").Value = 6277.6665
$ws.Range("K99").Value = 3963.4614
$ws.Range("L99").Value = 6277.6665
$ws.Range("M99").Value = -2465.4614
$ws.Range("N99").Value = -9273.666499999999

$ws.Range("H105").Value = 1179.5416
$ws.Range("I105").Value = 1149.909
$ws.Range("K105").Value = 1149.909
$ws.Range("M105").Value = 597.0909999999999

$ws.Range("H107").Value = 2256.1
$ws.Range("J107").Value = 1999.5
$ws.Range("L107").Value = 1999.5
$ws.Range("N107").Value = -5839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5659.8

$ws.Range("H31").Value = 25004164
$ws.Range("I31").Value = 38463930
$ws.Range("K31").Value = 38463930
$ws.Range("M31").Value = -38463635

$ws.Range("H34").Value = 25004164
$ws.Range("I34").Value = 38463930
$ws.Range("K34").Value = 38463930
$ws.Range("M34").Value = -38463728

$ws.Range("H113").Value = 5659.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 289
$ws.Range("I113").Value = 289
$ws.Range("K113").Value = 867
$ws.Range("M113").Value = 1303

$ws.Range("H121").Value = 59783.883
$ws.Range("I121").Value = 362.22223
$ws.Range("J121").Value = 126633.25
$ws.Range("K121").Value = 1086.66669
$ws.Range("L121").Value = 379899.75
$ws.Range("M121").Value = 223.33331
$ws.Range("N121").Value = -382519.75

$ws.Range("H122").Value = 6808.9
$ws.Range("J122").Value = 14295.667
$ws.Range("L122").Value = 128661.003
$ws.Range("N122").Value = -133561.003

$ws.Range("H140").Value = 3168.9285
$ws.Range("I140").Value = 2942.2727
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 8826.8181
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = -3646.8181
$ws.Range("N140").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2660.1035
$ws.Range("I80").Value = 2532.9285
$ws.Range("J80").Value = 2778.8
$ws.Range("K80").Value = 2532.9285
$ws.Range("L80").Value = 2778.8
$ws.Range("M80").Value = -1534.9285
$ws.Range("N80").Value = -4774.8

$ws.Range("H83").Value = 2660.1035
$ws.Range("I83").Value = 2532.9285
$ws.Range("J83").Value = 2778.8
$ws.Range("K83").Value = 12664.6425
$ws.Range("L83").Value = 13894
$ws.Range("M83").Value = -7672.6425
$ws.Range("N83").Value = -23878

$ws.Range("H125").Value = 25000
$ws.Range("J125").Value = 25000
$ws.Range("L125").Value = 25000
$ws.Range("N125").Value = -29920

$ws.Range("H132").Value = 5584.227
$ws.Range("I132").Value = 5068.1113
$ws.Range("K132").Value = 15204.3339
$ws.Range("M132").Value = -12674.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 271.8
$ws.Range("I16").Value = 271.8
$ws.Range("K16").Value = 271.8
$ws.Range("M16").Value = -101.8

$ws.Range("H40").Value = 3402.16
$ws.Range("I40").Value = 2178.5715
$ws.Range("K40").Value = 2178.5715
$ws.Range("M40").Value = -2042.5715

$ws.Range("H68").Value = 3437.25
$ws.Range("I68").Value = 2566.3333
$ws.Range("J68").Value = 6050
$ws.Range("K68").Value = 2566.3333
$ws.Range("L68").Value = 6050
$ws.Range("M68").Value = -1817.3333
$ws.Range("N68").Value = -7548

$ws.Range("H71").Value = 3437.25
$ws.Range("I71").Value = 2566.3333
$ws.Range("J71").Value = 6050
$ws.Range("K71").Value = 12831.6665
$ws.Range("L71").Value = 30250
$ws.Range("M71").Value = -9087.666499999999
$ws.Range("N71").Value = -37738

$ws.Range("H100").Value = 4645.4
$ws.Range("I100").Value = 2446.8333
$ws.Range("K100").Value = 2446.8333
$ws.Range("M100").Value = -1905.8333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1844222.2
$ws.Range("I4").Value = 46904.145
$ws.Range("K4").Value = 46904.145
$ws.Range("M4").Value = -46791.145

$ws.Range("H15").Value = 7495.25
$ws.Range("J15").Value = 7495.25
$ws.Range("L15").Value = 7495.25
$ws.Range("N15").Value = -8071.25

$ws.Range("H46").Value = 70969.25
$ws.Range("J46").Value = 74663.336
$ws.Range("L46").Value = 74663.336
$ws.Range("N46").Value = -75125.336

$ws.Range("H81").Value = 2557.6667
$ws.Range("I81").Value = 2869.2
$ws.Range("K81").Value = 5738.4
$ws.Range("M81").Value = -4677.4

$ws.Range("H84").Value = 2557.6667
$ws.Range("I84").Value = 2869.2
$ws.Range("K84").Value = 28692
$ws.Range("M84").Value = -23388

$ws.Range("H100").Value = 1250.5555
$ws.Range("I100").Value = 1376
$ws.Range("K100").Value = 2752
$ws.Range("M100").Value = -2211

$ws.Range("H126").Value = 6328.48
$ws.Range("I126").Value = 4490.35
$ws.Range("K126").Value = 13471.05
$ws.Range("M126").Value = -11001.05

$ws.Range("H132").Value = 3739.3918
$ws.Range("I132").Value = 3269.3442
$ws.Range("K132").Value = 9808.0326
$ws.Range("M132").Value = -7278.0326

$ws.Range("H134").Value = 70969.25
$ws.Range("J134").Value = 74663.336
$ws.Range("L134").Value = 223990.008
$ws.Range("N134").Value = -229060.008
